$d = $word.ActiveDocument

# --- Helper: find a paragraph whose trimmed text equals $text ---
function Find-ParagraphByText($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    throw "Paragraph with text '$text' not found"
}

# 1) Strike-through the "Teleports over different maps" bullet.
$pTeleport = Find-ParagraphByText $d "Teleports over different maps"
$pTeleport.Range.Font.StrikeThrough = 1

# 2) Strike-through the "Save/Load" sub-bullet.
$pSaveLoad = Find-ParagraphByText $d "Save/Load"
$pSaveLoad.Range.Font.StrikeThrough = 1

# 3) The "Publishment" bullet is currently split into two runs ("Publishm"
#    and "ent") with the hidden _GoBack bookmark sitting between them.
#    Move the bookmark onto the trailing (last, empty) paragraph first so
#    the two text runs become adjacent, then merge them into one run that
#    reads "Publishment".
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $pLast.Range)

$pPublish = Find-ParagraphByText $d "Publishment"
$pubStart = $pPublish.Range.Start

# "Publishm" occupies the first 8 characters of the paragraph; rewrite it
# in place to the full word.
$rFirst = $d.Range($pubStart, $pubStart + 8)
$rFirst.Text = "Publishment"

# The old "ent" run (3 characters) now trails right after, at offset +11;
# remove it since its text is already included in the rewritten run above.
$tailStart = $pubStart + 11
$tailEnd = $tailStart + 3
$rTail = $d.Range($tailStart, $tailEnd)
$rTail.Delete()
